# New forms for the BMGF demo: add "visit" and "plot" form entries,
# mirroring the existing "graphExample" entry pattern on the
# "survey" and "choices" sheets.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# Sheet "survey": append two 3-row blocks after row 52 (the
# existing "graphExample" block), one for "visit" and one for
# "plot", copying the formatting of the graphExample rows
# (A50:B52) so the same cell styles (s="3" / s="4") are reused
# instead of minting new ones.
# ---------------------------------------------------------------
$survey = $wb.Worksheets.Item("survey")

# --- visit block: rows 53-55 ---
$survey.Range("A51").Copy()
$survey.Range("A54").PasteSpecial(-4122)   # xlPasteFormats

$survey.Range("B51").Copy()
$survey.Range("B54").PasteSpecial(-4122)

$survey.Range("A52").Copy()
$survey.Range("A55").PasteSpecial(-4122)

$survey.Range("B52").Copy()
$survey.Range("B55").PasteSpecial(-4122)

$survey.Rows("53").RowHeight = 17.5
$survey.Rows("54").RowHeight = 66
$survey.Rows("55").RowHeight = 17

$survey.Range("A53").Value2 = "visit"
$survey.Range("B54").Value2 = "''?' + opendatakit.getHashString('../tables/visit/forms/visit/',null)"
$survey.Range("E54").Value2 = "external_link"
$survey.Range("G54").Value2 = "Open form"
$survey.Range("C55").Value2 = "exit section"

# --- plot block: rows 56-58 ---
$survey.Range("A51").Copy()
$survey.Range("A57").PasteSpecial(-4122)

$survey.Range("B51").Copy()
$survey.Range("B57").PasteSpecial(-4122)

$survey.Range("A52").Copy()
$survey.Range("A58").PasteSpecial(-4122)

$survey.Range("B52").Copy()
$survey.Range("B58").PasteSpecial(-4122)

$survey.Rows("56").RowHeight = 17.5
$survey.Rows("57").RowHeight = 66
$survey.Rows("58").RowHeight = 17

$survey.Range("A56").Value2 = "plot"
$survey.Range("B57").Value2 = "''?' + opendatakit.getHashString('../tables/plot/forms/plot/',null)"
$survey.Range("E57").Value2 = "external_link"
$survey.Range("G57").Value2 = "Open form"
$survey.Range("C58").Value2 = "exit section"

$survey.Application.CutCopyMode = $false

# Selection ends on the last new cell, matching the authored file.
$survey.Range("B58").Select()

# ---------------------------------------------------------------
# Sheet "choices": append the "visit" and "plot" test_forms
# choice rows, copying the format of the existing "graphExample"
# choice row (row 16) so column A reuses the same style (s="3").
# ---------------------------------------------------------------
$choices = $wb.Worksheets.Item("choices")

$choices.Range("A16").Copy()
$choices.Range("A17").PasteSpecial(-4122)
$choices.Range("A16").Copy()
$choices.Range("A18").PasteSpecial(-4122)
$choices.Application.CutCopyMode = $false

$choices.Range("A17").Value2 = "test_forms"
$choices.Range("B17").Value2 = "visit"
$choices.Range("C17").Value2 = "Visit"

$choices.Range("A18").Value2 = "test_forms"
$choices.Range("B18").Value2 = "plot"
$choices.Range("C18").Value2 = "Plot"

# choices is the workbook's active/selected tab; finish here so it
# stays the tab Excel reports as selected.
$choices.Activate()
$choices.Range("B19").Select()

# ---------------------------------------------------------------
# Workbook window position (best effort - cosmetic only).
# ---------------------------------------------------------------
$win = $wb.Windows.Item(1)
$win.Left = 10560
$win.Top = 1920
